# StyleIncludeQuotePrefix.xlsx was re-generated after ClosedXML's
# "adjust-to-content" (AutoFit) column-sizing logic was refactored, which
# produced slightly narrower auto-fitted widths for columns A:D:
#
#   cols A:B (style 0, no quote prefix): 4.070625  -> 2.996339
#   col  C   (style 1, quote prefix)   : 10.350625 -> 9.282054
#   col  D   (style 1, quote prefix)   : 12.840625 -> 11.710625
#
# Excel's ColumnWidth COM property is quantized to whole pixels on the
# Calibri-11 "Maximum Digit Width" grid (MDW = 6pt/char on this host):
#     stored_width = (round(ColumnWidth * 6) + 5) / 6
# so the ColumnWidth values below are chosen as the closest achievable
# input that rounds to the intended stored width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A:B -> stored width 2.996339 (closest achievable: 3)
$ws.Columns("A:B").ColumnWidth = 2.1666666666666665

# Column C -> stored width 9.282054 (closest achievable: 9.333333333333334)
$ws.Columns("C:C").ColumnWidth = 8.5

# Column D -> stored width 11.710625 (closest achievable: 11.666666666666666)
$ws.Columns("D:D").ColumnWidth = 10.833333333333334
